# Update the "取得日時" (retrieved datetime) column on the ランサーズ sheet
# from 2026-01-28 06:33:29 to 2026-01-28 06:42:02 for the existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2026-01-28 06:42:02"

for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
